$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.748.25'
$ws.Cells.Item(2, 5).Value = '  -0.24%  '
$ws.Cells.Item(3, 4).Value = '2.291.95'
$ws.Cells.Item(3, 5).Value = '  -2.06%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '98.65'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.69%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '270.51'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.58%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.618'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -1.40%  '
$ws.Cells.Item(8, 5).Value = '  -0.11%  '
$ws.Cells.Item(9, 5).Value = '  -3.53%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '44.99'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.68%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0928'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.05%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.92'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -2.92%  '
$ws.Cells.Item(13, 5).Value = '  +1.54%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.69'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.01%  '
$ws.Cells.Item(15, 4).Value = '2.635.28'
$ws.Cells.Item(15, 5).Value = '  -2.04%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.847'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -2.04%  '
$ws.Cells.Item(17, 4).Value = '2.293.77'
$ws.Cells.Item(17, 5).Value = '  -1.97%  '
$ws.Cells.Item(18, 4).Value = '43.734.82'
$ws.Cells.Item(18, 5).Value = '  -0.12%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000111'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.91%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.21'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -3.88%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '72.20'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.91%  '
$ws.Cells.Item(22, 5).Value = '  +8.46%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '232.92'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -3.03%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.86'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +12.78%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.09'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -3.66%  '
$ws.Cells.Item(26, 5).Value = '  +0.00%  '
$ws.Cells.Item(27, 5).Value = '  -1.87%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.24'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -1.68%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '38.36'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.64%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '176.75'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.63%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '21.81'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.50%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0891'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.38%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.43'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.23%  '
$ws.Cells.Item(35, 5).Value = '  +0.67%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.74'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +7.14%  '
$ws.Cells.Item(37, 5).Value = '  -0.55%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0351'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -2.63%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.52'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +3.64%  '
$ws.Cells.Item(40, 5).Value = '  -0.99%  '
$ws.Cells.Item(41, 5).Value = '  -2.68%  '
$ws.Cells.Item(42, 5).Value = '  -0.99%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.20'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.09%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.77'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +3.35%  '
$ws.Cells.Item(45, 5).Value = '  -4.63%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.22'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -3.20%  '
$ws.Cells.Item(47, 5).Value = '  -1.89%  '
$ws.Cells.Item(48, 5).Value = '  +0.96%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '98.19'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -2.73%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.55'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +13.01%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.440'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +5.23%  '
